# Implements feedback from the appellate efiling program review:
#   1. Remove the stray "_GoBack" bookmark that had been left after the
#      title's "Exemption from e-" (before "filing").
#   2. Fix the capitalization of "Circuit Clerk" -> "circuit clerk" in the
#      "File your form along with..." list item.
#   3. Fix the capitalization of "the Circuit Clerk." -> "the circuit
#      clerk." in the "Do not file these instructions..." paragraph. This
#      is where editing finished, so "_GoBack" now marks that spot
#      instead of its old location near the title.

$d = $word.ActiveDocument

# --- Step 1: drop the old "_GoBack" bookmark, wherever it currently sits ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete() | Out-Null
}

# --- Step 2: "... with the Circuit Clerk" -> "... with the circuit clerk" ---
# Locate the whole list-item sentence so the replacement XML below can
# rebuild every run in that paragraph (only the case of "Circuit Clerk"
# actually changes; the run boundaries mirror how Word split the edited
# text as it was retyped).
$r2 = $d.Content
$found2 = $r2.Find.Execute("File your form along with your other court papers with the Circuit Clerk.")
if (-not $found2) {
    throw "Could not find the 'File your form ... Circuit Clerk.' list item paragraph"
}
$target2 = $d.Range($r2.Start, $r2.End)

$xml2 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FD6B95" w:rsidRPr="00FD6B95" w:rsidRDefault="003E3F06" w:rsidP="00C812B0"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>File your</w:t></w:r><w:r w:rsidR="00FD6B95" w:rsidRPr="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>form</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> along with yo</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>ur other court papers with the circuit c</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>lerk</w:t></w:r><w:r w:rsidR="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target2.InsertXML($xml2) | Out-Null

# --- Step 3: "the Circuit Clerk." -> "the circuit clerk." --------------
# The edit (and thus the relocated "_GoBack" bookmark) lands between the
# retyped "the circuit c" and "lerk." pieces of the final bold run.
$r3 = $d.Content
$found3 = $r3.Find.Execute("Do not file these instructions with the Circuit Clerk.")
if (-not $found3) {
    throw "Could not find the 'Do not file ... Circuit Clerk.' paragraph"
}
$target3 = $d.Range($r3.Start, $r3.End)

$xml3 = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00FD6B95" w:rsidRPr="00FD6B95" w:rsidRDefault="00FD6B95" w:rsidP="00FD6B95"><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Do not </w:t></w:r><w:r w:rsidR="00C64446"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>file</w:t></w:r><w:r w:rsidRPr="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> these instructions </w:t></w:r><w:r w:rsidR="00C64446"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>with</w:t></w:r><w:r w:rsidRPr="00FD6B95"><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>the circuit c</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:rPr><w:rFonts w:ascii="Muli" w:hAnsi="Muli"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>lerk.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$target3.InsertXML($xml3) | Out-Null
